$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 label change (value C2 stays the same)
$ws.Range("B2").Value = "<then>"

# Row 3
$ws.Range("C3").Value = 6

# Row 4
$ws.Range("C4").Value = 8

# Row 5
$ws.Range("C5").Value = 12

# Row 6
$ws.Range("C6").Value = 9

# Row 7
$ws.Range("C7").Value = 5

# Row 8
$ws.Range("C8").Value = 8

# Row 9
$ws.Range("C9").Value = 4

# Row 10
$ws.Range("C10").Value = 10

# Row 13
$ws.Range("C13").Value = 11

# Row 15
$ws.Range("C15").Value = 7

# Row 16
$ws.Range("C16").Value = 5

# Row 17: B17 label change and C17 value change
$ws.Range("B17").Value = "<would>"
$ws.Range("C17").Value = 6
